# New crime data collected — weekly CompStat update for the 106th Precinct.
# Updates the report header (volume/number + week-covering dates) and
# refreshes the crime-complaint statistics table (rows 15-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates -------------------------------------------------
# "Volume 32   Number  10" -> "Volume 32   Number  11"
$ws.Range("A8").Characters(21, 2).Text = "11"

# "Report Covering the Week  3/3/2025  Through  3/9/2025"
#   -> "Report Covering the Week  3/10/2025  Through  3/16/2025"
$ws.Range("C9").Characters(27, 8).Text = "3/10/2025"
$ws.Range("C9").Characters(47, 8).Text = "3/16/2025"

# --- Row 15 (Rape) ---------------------------------------------------------
$ws.Range("D15").Value = "0"
$ws.Range("E15").Value = "***.*"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -50

# --- Row 16 (Robbery) -------------------------------------------------------
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = -41.176470588235
$ws.Range("I16").Value = 36
$ws.Range("J16").Value = 46
$ws.Range("K16").Value = -21.739130434782
$ws.Range("L16").Value = -20
$ws.Range("M16").Value = -40.983606557377
$ws.Range("N16").Value = -82.439024390243

# --- Row 17 (Fel. Assault) --------------------------------------------------
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -24
$ws.Range("I17").Value = 64
$ws.Range("J17").Value = 61
$ws.Range("K17").Value = 4.918032786885
$ws.Range("L17").Value = 8.474576271186
$ws.Range("M17").Value = 156
$ws.Range("N17").Value = -5.882352941176

# --- Row 18 (Burglary) ------------------------------------------------------
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 10
$ws.Range("I18").Value = 25
$ws.Range("J18").Value = 26
$ws.Range("K18").Value = -3.846153846153
$ws.Range("L18").Value = 4.166666666666
$ws.Range("M18").Value = -55.357142857142
$ws.Range("N18").Value = -90.909090909090

# --- Row 19 (Gr. Larceny) ---------------------------------------------------
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 47
$ws.Range("H19").Value = -31.914893617021
$ws.Range("I19").Value = 89
$ws.Range("J19").Value = 114
$ws.Range("K19").Value = -21.929824561403
$ws.Range("L19").Value = -23.931623931623
$ws.Range("M19").Value = 48.333333333333
$ws.Range("N19").Value = -25.833333333333

# --- Row 20 (G.L.A.) ---------------------------------------------------------
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 60
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = 76.470588235294
$ws.Range("I20").Value = 60
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = -14.285714285714
$ws.Range("L20").Value = 3.448275862068
$ws.Range("M20").Value = 11.111111111111
$ws.Range("N20").Value = -91.608391608391

# --- Row 21 (TOTAL) ----------------------------------------------------------
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -16.129032258064
$ws.Range("F21").Value = 104
$ws.Range("G21").Value = 117
$ws.Range("H21").Value = -11.111111111111
$ws.Range("I21").Value = 279
$ws.Range("J21").Value = 321
$ws.Range("K21").Value = -13.084112149532
$ws.Range("L21").Value = -9.120521172638
$ws.Range("M21").Value = 5.681818181818
$ws.Range("N21").Value = -80.028632784538

# --- Row 22 (Transit) ---------------------------------------------------------
$ws.Range("I22").Value = 1
$ws.Range("K22").Value = -80
$ws.Range("L22").Value = -80
$ws.Range("M22").Value = -80

# --- Row 24 (Petit Larceny) ---------------------------------------------------
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 4.545454545454
$ws.Range("F24").Value = 94
$ws.Range("G24").Value = 97
$ws.Range("H24").Value = -3.092783505154
$ws.Range("I24").Value = 263
$ws.Range("J24").Value = 240
$ws.Range("K24").Value = 9.583333333333
$ws.Range("L24").Value = 3.137254901960
$ws.Range("M24").Value = 97.744360902255

# --- Row 25 (Retail Theft) -----------------------------------------------------
$ws.Range("C25").Value = 17
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 52
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = 8.333333333333
$ws.Range("I25").Value = 144
$ws.Range("J25").Value = 113
$ws.Range("K25").Value = 27.433628318584
$ws.Range("L25").Value = 30.909090909090

# --- Row 26 (Misd. Assault) -----------------------------------------------------
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = 9.090909090909
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 111
$ws.Range("J26").Value = 107
$ws.Range("K26").Value = 3.738317757009
$ws.Range("L26").Value = 8.823529411764
$ws.Range("M26").Value = 21.978021978022

# --- Row 27 (UCR Rape*) -----------------------------------------------------------
$ws.Range("D27").Value = 1
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -16.666666666666
$ws.Range("L27").Value = -28.571428571428

# --- Row 28 (Other Sex Crimes) -----------------------------------------------------
$ws.Range("G28").Value = 2
$ws.Range("H28").Value = 150
$ws.Range("I28").Value = 12
$ws.Range("K28").Value = 33.333333333333
$ws.Range("L28").Value = -25
